$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (below the header row).
$ws.Rows.Item(2).Resize(2).Insert()

# The insert picks up the header row's formatting; clear it back to the
# plain (unstyled) look used by the rest of the data rows.
$ws.Rows.Item(2).Resize(2).ClearFormats()

# New row: Black Rock
$ws.Range("A2").Value = "Black Rock"
$ws.Range("B2").Value = "Woolworths Metro  40 Bluff Road, Black Rock VIC 3193"
$ws.Range("C2").Value = "30/12/20 5:30pm-5:55pm"
$ws.Range("D2").Value = "Case shopped"

# New row: Box Hill South
$ws.Range("A3").Value = "Box Hill South"
$ws.Range("B3").Value = "Bunnings  259 Middleborough Road, Box Hill South VIC 3128"
$ws.Range("C3").Value = "30/12/20 12:00pm-12:40pm"
$ws.Range("D3").Value = "Case shopped"

# Fix the Cheltenham exposure period (now at row 5 after the insert)
$ws.Range("C5").Value = "29/12/20 01:30pm-01:45pm"
